# "minor corrections to release"
#  - Drop the redundant "to gt" sheet (duplicate of "to gt (2)")
#  - Rename "to gt (2)" -> "table 2"
#  - Reset the saved scroll position on the "rounded" sheet so it opens at A1

$wb = $excel.ActiveWorkbook

# Remember which sheet is currently active so we can restore it afterwards.
$originalActiveSheetName = $wb.ActiveSheet.Name

# Remove the extra "to gt" worksheet entirely.
$wb.Worksheets.Item("to gt").Delete()

# Rename the remaining "to gt (2)" sheet to "table 2".
$wb.Worksheets.Item("to gt (2)").Name = "table 2"

# Clear the saved scroll position on the "rounded" sheet so it opens at A1
# instead of scrolled to Y5.
$roundedSheet = $wb.Worksheets.Item("rounded")
$roundedSheet.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1

# Restore the originally active tab (renamed if it was the "to gt (2)" sheet).
if ($originalActiveSheetName -eq "to gt (2)") {
    $wb.Worksheets.Item("table 2").Activate()
} else {
    $wb.Worksheets.Item($originalActiveSheetName).Activate()
}
